# Add files via upload
# - new column G with header "ok", plus two data cells ("sfp 600" / "sfp  ")
# - widen columns C and D
# - shrink a few row heights (5, 6, 7, 11) to match the new layout
# - adjust the final selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G ------------------------------------------------------
# G1: header cell - same look & feel as the other header cells (blue fill,
# bold white centered text) but paste-format from A1 first, then trim the
# border down to just a right edge (the look the real file ends up with).
$ws.Range("G1").Value = "ok"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Borders.Item(7).LineStyle = -4142
$ws.Range("G1").Borders.Item(8).LineStyle = -4142
$ws.Range("G1").Borders.Item(9).LineStyle = -4142
$ws.Range("G1").Borders.Item(10).Weight = -4138
$excel.CutCopyMode = $false

# G6 / G7: plain data cells
$ws.Range("G6").Value = "sfp 600"
$ws.Range("G7").Value = "sfp  "

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 29.33
$ws.Columns.Item(4).ColumnWidth = 34.17

# --- Row heights (rows got shorter once C/D widened and wrapped less) ----
$ws.Rows.Item(5).RowHeight = 30.75
$ws.Rows.Item(6).RowHeight = 45.75
$ws.Rows.Item(7).RowHeight = 30.75
$ws.Rows.Item(11).RowHeight = 30.75

# --- Final selection / view ----------------------------------------------
$ws.Range("G7").Select()
$excel.ActiveWindow.ScrollRow = 4
